# Scheduled-runner update: refresh currentAveragePrice / Leve profit figures
# across the per-job Sheets (market-board snapshot refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2579.3
$ws.Range("I113").Value = 3184.9285
$ws.Range("J113").Value = 2049.375
$ws.Range("K113").Value = 3184.9285
$ws.Range("L113").Value = 2049.375
$ws.Range("M113").Value = 69.07150000000001
$ws.Range("N113").Value = -8557.375

$ws.Range("H121").Value = 1155
$ws.Range("J121").Value = 1498.5
$ws.Range("L121").Value = 4495.5
$ws.Range("N121").Value = -7989.5

$ws.Range("H123").Value = 30620.6
$ws.Range("J123").Value = 30620.6
$ws.Range("L123").Value = 30620.6
$ws.Range("N123").Value = -40420.6

$ws.Range("H129").Value = 1225.6731
$ws.Range("J129").Value = 1158.0233
$ws.Range("L129").Value = 3474.0699
$ws.Range("N129").Value = -13474.0699

$ws.Range("H130").Value = 49085
$ws.Range("J130").Value = 49085
$ws.Range("L130").Value = 49085
$ws.Range("N130").Value = -59125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1480.4546
$ws.Range("I45").Value = 1476.4286
$ws.Range("J45").Value = 1487.5
$ws.Range("K45").Value = 1476.4286
$ws.Range("L45").Value = 1487.5
$ws.Range("M45").Value = -1099.4286
$ws.Range("N45").Value = -2241.5

$ws.Range("H122").Value = 1984.814
$ws.Range("I122").Value = 2053.1292
$ws.Range("J122").Value = 1808.3334
$ws.Range("K122").Value = 6159.3876
$ws.Range("L122").Value = 5425.0002
$ws.Range("M122").Value = -3709.3876
$ws.Range("N122").Value = -10325.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 685
$ws.Range("I94").Value = 620
$ws.Range("K94").Value = 620
$ws.Range("M94").Value = -169

$ws.Range("H105").Value = 2329.261
$ws.Range("I105").Value = 2039.8235
$ws.Range("K105").Value = 2039.8235
$ws.Range("M105").Value = -292.8235

$ws.Range("H107").Value = 2808.5293
$ws.Range("I107").Value = 2552.625
$ws.Range("J107").Value = 3036
$ws.Range("K107").Value = 2552.625
$ws.Range("L107").Value = 3036
$ws.Range("M107").Value = -632.625
$ws.Range("N107").Value = -6876

$ws.Range("H112").Value = 40795.2
$ws.Range("J112").Value = 40795.2
$ws.Range("L112").Value = 40795.2
$ws.Range("N112").Value = -43749.2

$ws.Range("H134").Value = 2781.99
$ws.Range("I134").Value = 1050.2683
$ws.Range("J134").Value = 3985.39
$ws.Range("K134").Value = 3150.8049
$ws.Range("L134").Value = 11956.17
$ws.Range("M134").Value = -615.8049000000001
$ws.Range("N134").Value = -17026.17

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2102
$ws.Range("I99").Value = 2162.4
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 2162.4
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -664.4000000000001
$ws.Range("N99").Value = -4796

$ws.Range("H126").Value = 2102
$ws.Range("I126").Value = 2162.4
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 6487.200000000001
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -4017.200000000001
$ws.Range("N126").Value = -10340

$ws.Range("H134").Value = 584952
$ws.Range("I134").Value = 1092.0625
$ws.Range("J134").Value = 1752671.9
$ws.Range("K134").Value = 3276.1875
$ws.Range("L134").Value = 5258015.699999999
$ws.Range("M134").Value = -741.1875
$ws.Range("N134").Value = -5263085.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 999999
$ws.Range("J105").Value = 999999
$ws.Range("L105").Value = 2999997
$ws.Range("N105").Value = -3005239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H102").Value = 1721.65
$ws.Range("I102").Value = 1338.6875
$ws.Range("J102").Value = 3253.5
$ws.Range("K102").Value = 1338.6875
$ws.Range("L102").Value = 3253.5
$ws.Range("M102").Value = 283.3125
$ws.Range("N102").Value = -6497.5

$ws.Range("H110").Value = 38651.832
$ws.Range("J110").Value = 38651.832
$ws.Range("L110").Value = 38651.832
$ws.Range("N110").Value = -46831.832

$ws.Range("H122").Value = 2001.1666
$ws.Range("I122").Value = 2176.75
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 6530.25
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -4080.25
$ws.Range("N122").Value = -9850

$ws.Range("H123").Value = 18914.666
$ws.Range("J123").Value = 18914.666
$ws.Range("L123").Value = 18914.666
$ws.Range("N123").Value = -23814.666

$ws.Range("H126").Value = 28575984
$ws.Range("I126").Value = 45460704
$ws.Range("J126").Value = 1841.3846
$ws.Range("K126").Value = 136382112
$ws.Range("L126").Value = 5524.1538
$ws.Range("M126").Value = -136379642
$ws.Range("N126").Value = -10464.1538

$ws.Range("H130").Value = 46611.8
$ws.Range("J130").Value = 46611.8
$ws.Range("L130").Value = 46611.8
$ws.Range("N130").Value = -56651.8

$ws.Range("H132").Value = 3013.9546
$ws.Range("I132").Value = 1231.6
$ws.Range("J132").Value = 6833.2856
$ws.Range("K132").Value = 3694.8
$ws.Range("L132").Value = 20499.8568
$ws.Range("M132").Value = -1164.8
$ws.Range("N132").Value = -25559.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3351.75
$ws.Range("I7").Value = 2773.7144
$ws.Range("J7").Value = 4161
$ws.Range("K7").Value = 2773.7144
$ws.Range("L7").Value = 4161
$ws.Range("M7").Value = -2661.7144
$ws.Range("N7").Value = -4385

$ws.Range("H40").Value = 2660.6
$ws.Range("I40").Value = 2601
$ws.Range("J40").Value = 2750
$ws.Range("K40").Value = 2601
$ws.Range("L40").Value = 2750
$ws.Range("M40").Value = -2465
$ws.Range("N40").Value = -3022

$ws.Range("H93").Value = 3186.2856
$ws.Range("I93").Value = 5000
$ws.Range("J93").Value = 2884
$ws.Range("K93").Value = 5000
$ws.Range("L93").Value = 2884
$ws.Range("M93").Value = -3752
$ws.Range("N93").Value = -5380

$ws.Range("H106").Value = 36689.332
$ws.Range("J106").Value = 36689.332
$ws.Range("L106").Value = 36689.332
$ws.Range("N106").Value = -39213.332

$ws.Range("H109").Value = 35277
$ws.Range("J109").Value = 35277
$ws.Range("L109").Value = 35277
$ws.Range("N109").Value = -38051

$ws.Range("H121").Value = 32713
$ws.Range("J121").Value = 32713
$ws.Range("L121").Value = 32713
$ws.Range("N121").Value = -36207

$ws.Range("H122").Value = 2128.5715
$ws.Range("I122").Value = 2100
$ws.Range("K122").Value = 6300
$ws.Range("M122").Value = -3850

$ws.Range("H126").Value = 3351.75
$ws.Range("I126").Value = 2773.7144
$ws.Range("J126").Value = 4161
$ws.Range("K126").Value = 8321.143199999999
$ws.Range("L126").Value = 12483
$ws.Range("M126").Value = -5851.143199999999
$ws.Range("N126").Value = -17423

$ws.Range("H132").Value = 5335.273
$ws.Range("I132").Value = 5249.5
$ws.Range("J132").Value = 5367.4375
$ws.Range("K132").Value = 15748.5
$ws.Range("L132").Value = 16102.3125
$ws.Range("M132").Value = -13218.5
$ws.Range("N132").Value = -21162.3125

$ws.Range("H134").Value = 50199.6
$ws.Range("J134").Value = 50199.6
$ws.Range("L134").Value = 50199.6
$ws.Range("N134").Value = -60339.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 38534.715
$ws.Range("J103").Value = 38534.715
$ws.Range("L103").Value = 38534.715
$ws.Range("N103").Value = -40878.715

$ws.Range("H126").Value = 4203981
$ws.Range("I126").Value = 4203981
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12611943
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12609473
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2565.0435
$ws.Range("I132").Value = 2014.1428
$ws.Range("J132").Value = 3422
$ws.Range("K132").Value = 6042.428400000001
$ws.Range("L132").Value = 10266
$ws.Range("M132").Value = -3512.428400000001
$ws.Range("N132").Value = -15326
